$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the timestamp column (O) for all data rows (2-397)
# from "2022-12-31 12:54:29" to "2022-12-31 20:49:19"
$ws.Range("O2:O397").Value = "2022-12-31 20:49:19"
